$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally only had a "Hello, World" banner (A2) and a TODAY()
# stamp (A3). This adds a small template table below it:
#   row 4: "name" | "rank"           (bold header, boxed)
#   row 5: "{{this}}"                (loop-open marker)
#   row 6: "{{name}}" | "{{rank}}"   (boxed data-row placeholders)
#   row 7: "{{/this}}"               (loop-close marker)
$ws.Range("A4").Value = "name"
$ws.Range("B4").Value = "rank"
$ws.Range("A5").Value = "{{this}}"
$ws.Range("A6").Value = "{{name}}"
$ws.Range("B6").Value = "{{rank}}"
$ws.Range("A7").Value = "{{/this}}"

# Box the placeholder row and the header row, and bold the header row.
$ws.Range("A6:B6").Borders.LineStyle = 1
$ws.Range("A4:B4").Borders.LineStyle = 1
$ws.Range("A4:B4").Font.Bold = $true

# Selection moved from E7 to E3.
$ws.Range("E3").Select() | Out-Null
